$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 19 (shifts old rows 19..42 down to 20..43)
$ws.Rows(19).Insert()

# Populate the newly inserted row 19 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,R repeat the same boilerplate values used by every
# other row in this sheet (same market / category / classification).
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44671
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112040
$ws.Cells.Item(19, 7).Value = "Cilantro"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 200
$ws.Cells.Item(19, 11).Value = 550
$ws.Cells.Item(19, 12).Value = 600
$ws.Cells.Item(19, 13).Value = 575
$ws.Cells.Item(19, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(19, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(19, 16).Value = 575
$ws.Cells.Item(19, 17).Value = 1
$ws.Cells.Item(19, 18).Value = "Hortaliza"
